$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1494"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "11051651.35"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1053"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6898667.66"
$ws.Range("D6").Style = "Normal"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "314"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1023801.26"
$ws.Range("D9").Style = "Normal"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "528"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3745714.31"
$ws.Range("D11").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "260"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1709569.34"
$ws.Range("D12").Style = "Normal"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "701"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6489679.30"
$ws.Range("D17").Style = "Normal"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "235"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "725674.14"
$ws.Range("D21").Style = "Normal"

$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "465"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3439193.87"
$ws.Range("D23").Style = "Normal"

$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "210"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1371733.23"
$ws.Range("D24").Style = "Normal"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "459"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1424472.48"
$ws.Range("D32").Style = "Normal"

$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "868"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6520658.93"
$ws.Range("D34").Style = "Normal"

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "559"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3415066.84"
$ws.Range("D36").Style = "Normal"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "514"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1909389.17"
$ws.Range("D44").Style = "Normal"

$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "197"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1193196.38"
$ws.Range("D45").Style = "Normal"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "283"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1318641.02"
$ws.Range("D46").Style = "Normal"

$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "1144"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8256160.48"
$ws.Range("D51").Style = "Normal"

$ws.Range("C52").NumberFormat = "@"
$ws.Range("C52").Value = "795"
$ws.Range("C52").Style = "Normal"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = "5064035.93"
$ws.Range("D52").Style = "Normal"

$ws.Range("C59").NumberFormat = "@"
$ws.Range("C59").Value = "6818"
$ws.Range("C59").Style = "Normal"
$ws.Range("D59").NumberFormat = "@"
$ws.Range("D59").Value = "34924356.76"
$ws.Range("D59").Style = "Normal"

$ws.Range("C69").NumberFormat = "@"
$ws.Range("C69").Value = "249"
$ws.Range("C69").Style = "Normal"
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value = "739482.59"
$ws.Range("D69").Style = "Normal"

$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "297"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "1019658.30"
$ws.Range("D75").Style = "Normal"

$ws.Range("C76").NumberFormat = "@"
$ws.Range("C76").Value = "510"
$ws.Range("C76").Style = "Normal"
$ws.Range("D76").NumberFormat = "@"
$ws.Range("D76").Value = "3194124.30"
$ws.Range("D76").Style = "Normal"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "304"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "2180416.07"
$ws.Range("D77").Style = "Normal"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "452"
$ws.Range("C80").Style = "Normal"
$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "1440772.96"
$ws.Range("D80").Style = "Normal"

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "1265"
$ws.Range("C82").Style = "Normal"
$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "9868167.91"
$ws.Range("D82").Style = "Normal"

$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "664"
$ws.Range("C83").Style = "Normal"
$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "4416787.10"
$ws.Range("D83").Style = "Normal"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "994"
$ws.Range("C86").Style = "Normal"
$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "2812211.79"
$ws.Range("D86").Style = "Normal"

$ws.Range("C89").NumberFormat = "@"
$ws.Range("C89").Value = "1363"
$ws.Range("C89").Style = "Normal"
$ws.Range("D89").NumberFormat = "@"
$ws.Range("D89").Value = "8506487.80"
$ws.Range("D89").Style = "Normal"

$ws.Range("C91").NumberFormat = "@"
$ws.Range("C91").Value = "969"
$ws.Range("C91").Style = "Normal"
$ws.Range("D91").NumberFormat = "@"
$ws.Range("D91").Value = "5402275.54"
$ws.Range("D91").Style = "Normal"

$ws.Range("C104").NumberFormat = "@"
$ws.Range("C104").Value = "1688"
$ws.Range("C104").Style = "Normal"
$ws.Range("D104").NumberFormat = "@"
$ws.Range("D104").Value = "9297300.59"
$ws.Range("D104").Style = "Normal"

$ws.Range("C106").NumberFormat = "@"
$ws.Range("C106").Value = "1651"
$ws.Range("C106").Style = "Normal"
$ws.Range("D106").NumberFormat = "@"
$ws.Range("D106").Value = "8630207.85"
$ws.Range("D106").Style = "Normal"

$ws.Range("C108").NumberFormat = "@"
$ws.Range("C108").Value = "82"
$ws.Range("C108").Style = "Normal"
$ws.Range("D108").NumberFormat = "@"
$ws.Range("D108").Value = "378788.23"
$ws.Range("D108").Style = "Normal"
